$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 154 — this shifts the existing rows 154:257 down to
# 155:258 (matching the diff, where every record from the old row N now
# lives at row N+1, and the workbook dimension grows from A1:R257 to
# A1:R258).
$ws.Rows("154:154").Insert()

# Populate the newly inserted row 154 with the new record's data.
$ws.Range("A154").Value2 = 8
$ws.Range("B154").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C154").Value2 = "Coquimbo"
$ws.Range("D154").Value2 = 45216
$ws.Range("E154").Value2 = 4
$ws.Range("F154").Value2 = 100112044
$ws.Range("G154").Value2 = "Perejil"
$ws.Range("H154").Value2 = "Sin especificar"
$ws.Range("I154").Value2 = "Primera"
$ws.Range("J154").Value2 = 3320
$ws.Range("K154").Value2 = 1500
$ws.Range("L154").Value2 = 2000
$ws.Range("M154").Value2 = 1831
$ws.Range("N154").Value2 = "`$/atado 1 a 1,5 kilos"
$ws.Range("O154").Value2 = "Provincia del Elquí"
$ws.Range("P154").Value2 = 1221
$ws.Range("Q154").Value2 = 1.5
$ws.Range("R154").Value2 = "Hortaliza"
